$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.332.44'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.863.00'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.05'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7001'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.53%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07895'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3118'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.31'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07791'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.37%  '
$ws.Range("D12").Value = '1.877.62'
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.139'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.25'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6965'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.556'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008536'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").Value = '29.382.44'
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.44'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = '2.123.31'
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.97'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.571'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1535'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.959'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.49'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.71'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.588'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.288'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.235'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.204'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05248'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.883'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7544'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.179'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.705'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("D38").Value = '1.276.17'
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01864'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.747'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8981'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.74'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.954'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.07'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.75%  '
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '2.022.33'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.584'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.791'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5173'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4281'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.28%  '
